$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 391
    3  = 393
    4  = 395
    5  = 396
    6  = 399
    7  = 400
    8  = 402
    9  = 405
    10 = 407
    11 = 409
    12 = 411
    13 = 413
    14 = 20
    15 = 26
    16 = 128
    17 = 138
    18 = 167
    19 = 196
    20 = 226
    21 = 250
    22 = 273
    23 = 293
    24 = 368
    25 = 418
    26 = 467
    27 = 489
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 1).Value = $values[$row]
}
